# "Changed the test runner path for features" - update the Doctors sheet
# sample rows from Delhi-based Internal Medicine doctors to Chennai-based
# Ear-Nose-Throat (ENT) Specialists.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Doctors")

$ws.Range("A2").Value = "Dr. Ramesh D"
$ws.Range("B2").Value = "Ear-Nose-Throat (ENT) Specialist"
$ws.Range("C2").Value = "31 years experience overall"
$ws.Range("D2").Value = "Ashok Nagar,Chennai"

$ws.Range("A3").Value = "Dr. Vikram P S J"
$ws.Range("B3").Value = "Ear-Nose-Throat (ENT) Specialist"
$ws.Range("C3").Value = "13 years experience overall"
$ws.Range("D3").Value = "Greams Road,Chennai"

$ws.Range("A4").Value = "Dr. Balaji R"
$ws.Range("B4").Value = "Ear-Nose-Throat (ENT) Specialist"
$ws.Range("C4").Value = "24 years experience overall"
$ws.Range("D4").Value = "Thousand Lights,Chennai"

$ws.Range("A5").Value = "Dr. Nitya Narayanan"
$ws.Range("B5").Value = "Ear-Nose-Throat (ENT) Specialist"
$ws.Range("C5").Value = "29 years experience overall"
$ws.Range("D5").Value = "Thousand Lights,Chennai"

$ws.Range("A6").Value = "Dr. Sudha Anantha Krishnan"
$ws.Range("B6").Value = "Ear-Nose-Throat (ENT) Specialist"
$ws.Range("C6").Value = "27 years experience overall"
$ws.Range("D6").Value = "Kilpauk,Chennai"
